# Weekly driver report update for 2025-04-21
# Updates the "Bad Drivers" table on the "Driver Summary" sheet:
#  - refreshes Critical Minutes / Good Roaming % for several existing drivers
#  - inserts a newly-observed driver row
#  - re-sorts a couple of rows that now fall in a different rank order
#  - refreshes the Totals row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 4-6 with refreshed sample counts -----------------
$ws.Range("C4").Value = 3426
$ws.Range("D4").Value = 96

$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 97.8

$ws.Range("C6").Value = 1121
$ws.Range("D6").Value = 98

# --- Insert newly-observed driver row at row 7 -----------------------------
# This pushes the previous rows 7-52 down to 8-53 (matches new dimension A1:J53)
$ws.Rows(7).Insert()

$ws.Range("A7").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 98.40000000000001

# --- Row 8 (previously row 7) - refreshed Critical Minutes ----------------
$ws.Range("C8").Value = 7759

# --- Rows 10 & 11 now swap order (21.50.1.1 moves above 22.80.1.1) ---------
$ws.Range("A10").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.50.1.1"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 23
$ws.Range("D10").Value = 98.8

$ws.Range("A11").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.80.1.1"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 98.8

# --- Totals row (now row 12) -----------------------------------------------
$ws.Range("B12").Value = 124
$ws.Range("C12").Value = 12348
